$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (shifts boson..syst6_c from H:W to I:X)
$ws.Columns("H").Insert()

# Header for the newly inserted column
$ws.Range("H1").Value = "pt_max"

# Fill the new column's data rows (2-11) with the constant 50
$ws.Range("H2:H11").Value = 50

# Re-assert the (now shifted) "syst_scaled"-derived formula column so the
# engine keeps it as one shared formula group, matching the other two
# shared formula columns that survived the column insert automatically.
$ws.Range("O2:O11").Formula = "=L2/100"

# Match the author's final selection state
$ws.Range("H16").Select() | Out-Null

Write-Output "done"
